# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with the latest scraped figures. The Price column must stay plain text
# (values like "42.774.69" or "0.628" should not turn into numbers), so we
# force the column to Text format before writing, then drop back to the
# default "Normal" style afterwards to avoid leaving an explicit format on
# the cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '42.774.69'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '2.299.83'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '317.83'
$ws.Range("E5").Value = '  -0.48%  '
$ws.Range("D6").Value = '104.42'
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  -0.61%  '
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").Value = '0.604'
$ws.Range("E9").Value = '  -1.38%  '
$ws.Range("D10").Value = '39.83'
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("D11").Value = '0.0909'
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("D12").Value = '8.54'
$ws.Range("E12").Value = '  +1.94%  '
$ws.Range("E13").Value = '  +2.31%  '
$ws.Range("E14").Value = '  +3.97%  '
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").Value = '2.652.83'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '2.303.10'
$ws.Range("E17").Value = '  -0.95%  '
$ws.Range("D18").Value = '42.719.10'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("D19").Value = '15.06'
$ws.Range("E19").Value = '  +37.43%  '
$ws.Range("D20").Value = '7.55'
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E21").Value = '  +0.17%  '
$ws.Range("E22").Value = '  +0.88%  '
$ws.Range("D23").Value = '3.56'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("D24").Value = '266.83'
$ws.Range("E24").Value = '  -5.70%  '
$ws.Range("D25").Value = '2.23'
$ws.Range("E25").Value = '  -2.40%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = '10.96'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = '6.80'
$ws.Range("E29").Value = '  +14.43%  '
$ws.Range("D30").Value = '22.63'
$ws.Range("E30").Value = '  -1.64%  '
$ws.Range("D31").Value = '37.43'
$ws.Range("E31").Value = '  +3.34%  '
$ws.Range("D32").Value = '165.94'
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").Value = '0.0883'
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("E34").Value = '  -4.45%  '
$ws.Range("D36").Value = '0.113'
$ws.Range("E36").Value = '  -3.45%  '
$ws.Range("D37").Value = '4.57'
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("E38").Value = '  -4.62%  '
$ws.Range("D39").Value = '3.74'
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").Value = '2.71'
$ws.Range("E40").Value = '  -3.24%  '
$ws.Range("D41").Value = '1.59'
$ws.Range("E41").Value = '  +4.90%  '
$ws.Range("D42").Value = '70.53'
$ws.Range("E42").Value = '  +0.75%  '
$ws.Range("D43").Value = '0.229'
$ws.Range("E43").Value = '  +0.65%  '
$ws.Range("D44").Value = '95.37'
$ws.Range("E44").Value = '  -3.72%  '
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '12.30'
$ws.Range("E46").Value = '  +0.71%  '
$ws.Range("D47").Value = '115.78'
$ws.Range("E47").Value = '  +2.38%  '
$ws.Range("D48").Value = '81.15'
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").Value = '1.696.24'
$ws.Range("E49").Value = '  +5.40%  '
$ws.Range("D50").Value = '8.82'
$ws.Range("E50").Value = '  -1.37%  '
$ws.Range("D51").Value = '5.17'
$ws.Range("E51").Value = '  -2.59%  '

$dRange.Style = "Normal"
